$d = $word.ActiveDocument

$titleText = "Play Choy Sun Doa Free: Review of Aristocrat's Exciting Slot"
$oldMetaText = "Get ready for an immersive and engaging online slot experience. Play Choy Sun Doa free and enjoy 243 ways to win and exciting bonus features."
$newPromptText = 'Create a feature image for "Choy Sun Doa" that features a happy Maya warrior with glasses. The image should be in a cartoon style and incorporate elements from the game, such as the God of Wealth symbol and the gold bullion. The warrior should be holding a winning ticket or surrounded by gold coins, showcasing the game''s potential for big payouts. The overall design should be bright, colorful, and fun, capturing the essence of the game''s positive and exciting atmosphere.'

# ------------------------------------------------------------------
# Part 1: right after the H1 title paragraph, insert a new paragraph
# holding the "Meta description: ..." line (bold label + plain text).
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' +
  '<w:r/>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
  '<w:r><w:t>: ' + $oldMetaText + '</w:t></w:r>' +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$null = $metaPara.Range.InsertXML($metaXml)

# ------------------------------------------------------------------
# Part 2: near the end of the document there used to be a leftover,
# duplicated copy of the title paragraph followed by the original
# meta-description paragraph (italic). Drop the duplicated title
# paragraph entirely, and turn the remaining italic paragraph into
# the feature-image generation prompt.
# ------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq $titleText) {
        $dupRange = $d.Range($p.Range.Start, $p.Range.End)
        $dupRange.Delete()
        break
    }
}

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq $oldMetaText) {
        $textRange = $d.Range($p.Range.Start, $p.Range.End - 1)
        $textRange.Text = $newPromptText
        break
    }
}
